$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.539087912495223
$ws.Range("J2").Value = 0.539087912495223
$ws.Range("M2").Value = 0.003710666666666666
$ws.Range("N2").Value = 0.011132
$ws.Range("O2").Value = 0.001642024256586498
$ws.Range("P2").Value = 0.001642024256586498
$ws.Range("Q2").Value = 0.0005216999431111111
$ws.Range("R2").Value = 0.004695299488
$ws.Range("S2").Value = 0.0008851954287497359
$ws.Range("T2").Value = 0.0008851954287497359

# Row 3
$ws.Range("I3").Value = 0.539087912495223
$ws.Range("J3").Value = 0.539087912495223
$ws.Range("M3").Value = 2.049608666666666
$ws.Range("N3").Value = 6.148826
$ws.Range("O3").Value = 0.9069818039462568
$ws.Range("P3").Value = 0.9069818039462569
$ws.Range("S3").Value = 0.4889429273605392
$ws.Range("T3").Value = 0.4889429273605393

# Row 4
$ws.Range("I4").Value = 0.539087912495223
$ws.Range("J4").Value = 0.539087912495223
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.206493
$ws.Range("N4").Value = 0.619479
$ws.Range("O4").Value = 0.09137617179715662
$ws.Range("P4").Value = 0.09137617179715661
$ws.Range("Q4").Value = 0.029031814504
$ws.Range("R4").Value = 0.261286330536
$ws.Range("S4").Value = 0.04925978970593404
$ws.Range("T4").Value = 0.04925978970593403

# Row 5
$ws.Range("G5").Value = 0.1202063333333333
$ws.Range("H5").Value = 0.360619
$ws.Range("I5").Value = 0.460912087504777
$ws.Range("J5").Value = 0.460912087504777
$ws.Range("M5").Value = 0.003710666666666666
$ws.Range("N5").Value = 0.011132
$ws.Range("O5").Value = 0.001642024256586498
$ws.Range("P5").Value = 0.001642024256586498
$ws.Range("Q5").Value = 0.0004460456342222222
$ws.Range("R5").Value = 0.004014410708
$ws.Range("S5").Value = 0.0007568288278367624
$ws.Range("T5").Value = 0.0007568288278367625

# Row 6
$ws.Range("G6").Value = 0.1202063333333333
$ws.Range("H6").Value = 0.360619
$ws.Range("I6").Value = 0.460912087504777
$ws.Range("J6").Value = 0.460912087504777
$ws.Range("M6").Value = 2.049608666666666
$ws.Range("N6").Value = 6.148826
$ws.Range("O6").Value = 0.9069818039462568
$ws.Range("P6").Value = 0.9069818039462569
$ws.Range("Q6").Value = 0.2463759425882222
$ws.Range("R6").Value = 2.217383483294
$ws.Range("S6").Value = 0.4180388765857176
$ws.Range("T6").Value = 0.4180388765857177

# Row 7
$ws.Range("G7").Value = 0.1202063333333333
$ws.Range("H7").Value = 0.360619
$ws.Range("I7").Value = 0.460912087504777
$ws.Range("J7").Value = 0.460912087504777
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.206493
$ws.Range("N7").Value = 0.619479
$ws.Range("O7").Value = 0.09137617179715662
$ws.Range("P7").Value = 0.09137617179715661
$ws.Range("Q7").Value = 0.024821766389
$ws.Range("R7").Value = 0.223395897501
$ws.Range("S7").Value = 0.04211638209122259
$ws.Range("T7").Value = 0.04211638209122259
